# Updated cryptos list (prices + 1h volume %) as produced by the
# GitHub Actions scraper run. Price-column cells whose new text looks
# like a plain number (e.g. "263.28") are forced to stay text - set
# NumberFormat to Text ("@") before assigning the string, then restore
# the "Normal" style so the cell keeps its original (unstyled) look -
# otherwise Excel's COM layer would silently reinterpret them as
# floating-point numbers, same as the workbook's existing inline-string
# cells such as "26.590.90" (not valid numbers, so unaffected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.565.48'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '1.851.46'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '263.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5257'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3238'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06811'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7843'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07771'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('D13').Value = '1.853.83'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.038'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.99'
$ws.Range('D17').Style = 'Normal'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007961'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '26.594.19'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.644'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.484'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.009'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.172'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.680'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '111.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08721'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.105'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7228'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.14%  '
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.875'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.104'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.272'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01793'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.4864'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9007'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.957'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.68%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.689'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4192'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05883'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.051'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1237'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.8912'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.56%  '
